# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Zafiro rojo / Zafiro verde, Primera,
# Limache, $/caja 18 kilos, fecha 44637) above the existing row 548,
# pushing the rest of the "Pimiento" block down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 548..646 down to 550..648, leaving 548:549 blank.
$ws.Rows("548:549").Insert()

# New row 548: Zafiro rojo / Primera
$ws.Range("A548").Value = 4
$ws.Range("B548").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C548").Value = "Los Lagos"
$ws.Range("D548").Value = 44637
$ws.Range("E548").Value = 10
$ws.Range("F548").Value = 100112002
$ws.Range("G548").Value = "Pimiento"
$ws.Range("H548").Value = "Zafiro rojo"
$ws.Range("I548").Value = "Primera"
$ws.Range("J548").Value = 120
$ws.Range("K548").Value = 31000
$ws.Range("L548").Value = 31000
$ws.Range("M548").Value = 31000
$ws.Range("N548").Value = "$/caja 18 kilos"
$ws.Range("O548").Value = "Limache"
$ws.Range("P548").Value = 1722
$ws.Range("Q548").Value = 18
$ws.Range("R548").Value = "Hortaliza"

# New row 549: Zafiro verde / Primera
$ws.Range("A549").Value = 4
$ws.Range("B549").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C549").Value = "Los Lagos"
$ws.Range("D549").Value = 44637
$ws.Range("E549").Value = 10
$ws.Range("F549").Value = 100112002
$ws.Range("G549").Value = "Pimiento"
$ws.Range("H549").Value = "Zafiro verde"
$ws.Range("I549").Value = "Primera"
$ws.Range("J549").Value = 120
$ws.Range("K549").Value = 22000
$ws.Range("L549").Value = 22000
$ws.Range("M549").Value = 22000
$ws.Range("N549").Value = "$/caja 18 kilos"
$ws.Range("O549").Value = "Limache"
$ws.Range("P549").Value = 1222
$ws.Range("Q549").Value = 18
$ws.Range("R549").Value = "Hortaliza"
